# Doris documentation sample: rename several header columns so that the
# "code" value of each cause-of-death is distinguished from its free-text
# label, and widen the CauseOfDeathURIA column a bit.
#
#   Age                -> EstimatedAge
#   CauseOfDeathA       -> CauseOfDeathCodeA
#   CauseOfDeathB       -> CauseOfDeathCodeB
#   CauseOfDeathC       -> CauseOfDeathCodeC
#   CauseOfDeathD       -> CauseOfDeathCodeD
#   CauseOfDeathE       -> CauseOfDeathCodeE
#   CauseOfDeathPart2  -> CauseOfDeathCodePart2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "EstimatedAge"
$ws.Range("H1").Value = "CauseOfDeathCodeA"
$ws.Range("L1").Value = "CauseOfDeathCodeB"
$ws.Range("P1").Value = "CauseOfDeathCodeC"
$ws.Range("T1").Value = "CauseOfDeathCodeD"
$ws.Range("X1").Value = "CauseOfDeathCodeE"
$ws.Range("AB1").Value = "CauseOfDeathCodePart2"

# Widen column I (CauseOfDeathURIA) slightly.
$ws.Columns.Item(9).ColumnWidth = 22.14
